$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift existing data rows 2-21 down by one row (to 3-22) ---
# Copy from the bottom up so we never overwrite a row before it has been
# read. Using .Value2 (not .Value) for accurate round-tripping of numbers.
$lastExistingRow = 21
for ($r = $lastExistingRow; $r -ge 2; $r--) {
    $ws.Cells.Item($r + 1, 1).Value2 = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r + 1, 2).Value2 = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r + 1, 3).Value2 = $ws.Cells.Item($r, 3).Value2
}

# --- Write the new row of data into the now-empty row 2 ---
$ws.Cells.Item(2, 1).Value2 = -0.09987647086381909
$ws.Cells.Item(2, 2).Value2 = -0.8263479471206665
$ws.Cells.Item(2, 3).Value2 = 0.256868839263916

# --- Append 9 brand-new rows of data after the (shifted) last row, 22 ---
$newRows = @(
    @(0.4167627990245819, 1.646892666816711, 0.836885392665863),
    @(0.5007568001747131, -1.671174645423889, -0.2092213481664657),
    @(0.0099265603348612, -3.577379703521729, -0.4952589869499206),
    @(-1.279915452003479, -6.508005619049072, -0.1085812970995903),
    @(-0.9758572578430176, -1.915215253829956, 1.272432327270508),
    @(0.0806342139840126, -1.140027284622192, 0.0485637858510017),
    @(0.8413141369819641, -0.3859141170978546, 0.5590944290161133),
    @(-0.1750128865242004, 2.383749008178711, -0.1401935666799545),
    @(0.4100432991981506, 2.355190992355347, 0.1059851199388504)
)

$startRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value2 = $newRows[$i][2]
}
